# Swap the "B quarter" row and "C quarter" row contents for each year
# 2003-2018. Row numbers stay fixed; the full row contents (A:O) are
# exchanged between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(9, 10),
    @(13, 14),
    @(17, 18),
    @(21, 22),
    @(25, 26),
    @(29, 30),
    @(33, 34),
    @(37, 38),
    @(41, 42),
    @(45, 46),
    @(49, 50),
    @(53, 54),
    @(57, 58),
    @(61, 62),
    @(65, 66),
    @(69, 70)
)

foreach ($pair in $rowPairs) {
    $rowB = $pair[0]
    $rowC = $pair[1]

    $rangeB = $ws.Range("A$rowB`:O$rowB")
    $rangeC = $ws.Range("A$rowC`:O$rowC")

    $valuesB = $rangeB.Value2
    $valuesC = $rangeC.Value2

    $rangeB.Value2 = $valuesC
    $rangeC.Value2 = $valuesB
}
